$d = $word.ActiveDocument

# 1. Update the cached/displayed "Reported" date field text.
$d.Content.Find.Execute("10-Jan-2023", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "9-Aug-2023", 2)

# 2. Remove the leftover "_GoBack" bookmark (last-edit marker) at the end
#    of the document body.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Shrink the CLIN2(SUBHEADINGS) style (and its linked character style)
#    from 13pt to 11pt.
$paraStyle = $d.Styles("CLIN2(SUBHEADINGS)")
$paraStyle.Font.Size = 11

$charStyle = $d.Styles("CLIN2(SUBHEADINGS) Char")
$charStyle.Font.Size = 11
